$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-05-17 Friday" "2024-05-18 Saturday"

Replace-Text "697÷5=139, 2" "924÷2=462, 0"
Replace-Text "520÷2=260, 0" "889÷2=444, 1"
Replace-Text "336÷8=42, 0" "370÷2=185, 0"
Replace-Text "778÷7=111, 1" "977÷7=139, 4"
Replace-Text "318÷6=53, 0" "470÷6=78, 2"

Replace-Text "854÷9=94, 8" "460÷3=153, 1"
Replace-Text "505÷6=84, 1" "111÷8=13, 7"
Replace-Text "133÷7=19, 0" "836÷2=418, 0"
Replace-Text "119÷5=23, 4" "409÷8=51, 1"
Replace-Text "699÷6=116, 3" "452÷8=56, 4"

Replace-Text "777÷7=111, 0" "351÷4=87, 3"
Replace-Text "140÷3=46, 2" "741÷6=123, 3"
Replace-Text "595÷7=85, 0" "381÷9=42, 3"
Replace-Text "224÷4=56, 0" "296÷9=32, 8"
Replace-Text "265÷5=53, 0" "669÷2=334, 1"

Replace-Text "623÷3=207, 2" "531÷2=265, 1"
Replace-Text "829÷4=207, 1" "335÷2=167, 1"
Replace-Text "145÷6=24, 1" "815÷8=101, 7"
Replace-Text "857÷9=95, 2" "486÷2=243, 0"
Replace-Text "831÷9=92, 3" "865÷8=108, 1"

Replace-Text "586÷8=73, 2" "381÷6=63, 3"
Replace-Text "772÷8=96, 4" "507÷2=253, 1"
Replace-Text "374÷2=187, 0" "679÷8=84, 7"
Replace-Text "790÷9=87, 7" "472÷5=94, 2"
Replace-Text "105÷8=13, 1" "676÷3=225, 1"
